$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header row: new column labels Wins / Losses / Ties
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Apply the same style as the other header cells (copy from AC1)
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) # xlPasteFormats

for ($r = 2; $r -le 55; $r++) {
    $ws.Cells.Item($r, 30).Value = 66   # AD
    $ws.Cells.Item($r, 31).Value = 96   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
